# Apply "realistic run params" changes to config.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the stray Debug block that lived at F1:H2 (duplicate of A31:D32)
$ws.Range("F1:H2").Clear()

# Rename "Print Plots?" (A10) to "Output Plots?" and update its note (D10)
$ws.Range("A10").Value = "Output Plots?"
$ws.Range("D10").Value = "If set to no, nothing in this section matters."

# Turn on Plot Contours / Plot Hatches
$ws.Range("B13").Value = "Yes"
$ws.Range("B14").Value = "Yes"

# Update selection / active cell to match the saved state
$ws.Range("D11").Select()

$wb.Save()
